$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3076.6667
$ws.Cells.Item(40, 10).Value = 3153.3333
$ws.Cells.Item(40, 12).Value = 3153.3333
$ws.Cells.Item(40, 14).Value = -3503.3333
$ws.Cells.Item(62, 8).Value = 3189.9
$ws.Cells.Item(62, 9).Value = 3189.9
$ws.Cells.Item(62, 11).Value = 3189.9
$ws.Cells.Item(62, 13).Value = -2565.9
$ws.Cells.Item(65, 8).Value = 3189.9
$ws.Cells.Item(65, 9).Value = 3189.9
$ws.Cells.Item(65, 11).Value = 15949.5
$ws.Cells.Item(65, 13).Value = -12829.5
$ws.Cells.Item(138, 8).Value = 1763.73
$ws.Cells.Item(138, 10).Value = 2142.3582
$ws.Cells.Item(138, 12).Value = 6427.0746
$ws.Cells.Item(138, 14).Value = -16707.0746

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 556695.4
$ws.Cells.Item(2, 9).Value = 926692.3
$ws.Cells.Item(2, 10).Value = 1700
$ws.Cells.Item(2, 11).Value = 926692.3
$ws.Cells.Item(2, 12).Value = 1700
$ws.Cells.Item(2, 13).Value = -926579.3
$ws.Cells.Item(2, 14).Value = -1926
$ws.Cells.Item(32, 8).Value = 4432.8774
$ws.Cells.Item(32, 9).Value = 2945.3372
$ws.Cells.Item(32, 11).Value = 2945.3372
$ws.Cells.Item(32, 13).Value = -2658.3372
$ws.Cells.Item(61, 8).Value = 30643.035
$ws.Cells.Item(61, 9).Value = 34400.207
$ws.Cells.Item(61, 11).Value = 34400.207
$ws.Cells.Item(61, 13).Value = -34188.207
$ws.Cells.Item(74, 8).Value = 660.12195
$ws.Cells.Item(74, 9).Value = 531.625
$ws.Cells.Item(74, 11).Value = 531.625
$ws.Cells.Item(74, 13).Value = 342.375
$ws.Cells.Item(77, 8).Value = 660.12195
$ws.Cells.Item(77, 9).Value = 531.625
$ws.Cells.Item(77, 11).Value = 2658.125
$ws.Cells.Item(77, 13).Value = 1709.875
$ws.Cells.Item(102, 8).Value = 1399.2
$ws.Cells.Item(102, 9).Value = 1399.2
$ws.Cells.Item(102, 11).Value = 1399.2
$ws.Cells.Item(102, 13).Value = 222.8
$ws.Cells.Item(116, 8).Value = 556695.4
$ws.Cells.Item(116, 9).Value = 926692.3
$ws.Cells.Item(116, 10).Value = 1700
$ws.Cells.Item(116, 11).Value = 926692.3
$ws.Cells.Item(116, 12).Value = 1700
$ws.Cells.Item(116, 13).Value = -924398.3
$ws.Cells.Item(116, 14).Value = -6288
$ws.Cells.Item(136, 8).Value = 30643.035
$ws.Cells.Item(136, 9).Value = 34400.207
$ws.Cells.Item(136, 11).Value = 103200.621
$ws.Cells.Item(136, 13).Value = -100650.621

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 556695.4
$ws.Cells.Item(3, 9).Value = 926692.3
$ws.Cells.Item(3, 10).Value = 1700
$ws.Cells.Item(3, 11).Value = 926692.3
$ws.Cells.Item(3, 12).Value = 1700
$ws.Cells.Item(3, 13).Value = -926578.3
$ws.Cells.Item(3, 14).Value = -1928
$ws.Cells.Item(20, 8).Value = 2812.4482
$ws.Cells.Item(20, 9).Value = 2518.9473
$ws.Cells.Item(20, 11).Value = 2518.9473
$ws.Cells.Item(20, 13).Value = -2271.9473
$ws.Cells.Item(105, 8).Value = 2096.2415
$ws.Cells.Item(105, 9).Value = 2180.12
$ws.Cells.Item(105, 11).Value = 2180.12
$ws.Cells.Item(105, 13).Value = -433.1199999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1421.6072
$ws.Cells.Item(31, 9).Value = 902.44446
$ws.Cells.Item(31, 10).Value = 2356.1
$ws.Cells.Item(31, 11).Value = 902.44446
$ws.Cells.Item(31, 12).Value = 2356.1
$ws.Cells.Item(31, 13).Value = -607.44446
$ws.Cells.Item(31, 14).Value = -2946.1
$ws.Cells.Item(34, 8).Value = 1421.6072
$ws.Cells.Item(34, 9).Value = 902.44446
$ws.Cells.Item(34, 10).Value = 2356.1
$ws.Cells.Item(34, 11).Value = 902.44446
$ws.Cells.Item(34, 12).Value = 2356.1
$ws.Cells.Item(34, 13).Value = -700.44446
$ws.Cells.Item(34, 14).Value = -2760.1
$ws.Cells.Item(51, 8).Value = 30000
$ws.Cells.Item(51, 10).Value = 35000
$ws.Cells.Item(51, 12).Value = 35000
$ws.Cells.Item(51, 14).Value = -36472
$ws.Cells.Item(61, 8).Value = 30000
$ws.Cells.Item(61, 10).Value = 35000
$ws.Cells.Item(61, 12).Value = 35000
$ws.Cells.Item(61, 14).Value = -35696
$ws.Cells.Item(105, 8).Value = 1356
$ws.Cells.Item(105, 9).Value = 1379.5
$ws.Cells.Item(105, 10).Value = 1238.5
$ws.Cells.Item(105, 11).Value = 1379.5
$ws.Cells.Item(105, 12).Value = 1238.5
$ws.Cells.Item(105, 13).Value = 367.5
$ws.Cells.Item(105, 14).Value = -4732.5
$ws.Cells.Item(132, 8).Value = 1583.6666
$ws.Cells.Item(132, 9).Value = 1186.4667
$ws.Cells.Item(132, 10).Value = 2245.6667
$ws.Cells.Item(132, 11).Value = 3559.4001
$ws.Cells.Item(132, 12).Value = 6737.000100000001
$ws.Cells.Item(132, 13).Value = -1029.4001
$ws.Cells.Item(132, 14).Value = -11797.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 250.55556
$ws.Cells.Item(2, 9).Value = 286.25
$ws.Cells.Item(2, 10).Value = 222
$ws.Cells.Item(2, 11).Value = 1717.5
$ws.Cells.Item(2, 12).Value = 1332
$ws.Cells.Item(2, 13).Value = -1604.5
$ws.Cells.Item(2, 14).Value = -1558
$ws.Cells.Item(56, 8).Value = 6380.727
$ws.Cells.Item(56, 9).Value = 6380.727
$ws.Cells.Item(56, 11).Value = 6380.727
$ws.Cells.Item(56, 13).Value = -5850.727
$ws.Cells.Item(69, 8).Value = 2575.4736
$ws.Cells.Item(69, 9).Value = 1899.6
$ws.Cells.Item(69, 11).Value = 5698.799999999999
$ws.Cells.Item(69, 13).Value = -4887.799999999999
$ws.Cells.Item(72, 8).Value = 2575.4736
$ws.Cells.Item(72, 9).Value = 1899.6
$ws.Cells.Item(72, 11).Value = 17096.4
$ws.Cells.Item(72, 13).Value = -13040.4
$ws.Cells.Item(80, 8).Value = 1850
$ws.Cells.Item(80, 9).Value = 1125
$ws.Cells.Item(80, 11).Value = 3375
$ws.Cells.Item(80, 13).Value = -2439
$ws.Cells.Item(83, 8).Value = 1850
$ws.Cells.Item(83, 9).Value = 1125
$ws.Cells.Item(83, 11).Value = 10125
$ws.Cells.Item(83, 13).Value = -5445
$ws.Cells.Item(131, 8).Value = 38293.316
$ws.Cells.Item(131, 10).Value = 46719.723
$ws.Cells.Item(131, 12).Value = 140159.169
$ws.Cells.Item(131, 14).Value = -150239.169

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3658.75
$ws.Cells.Item(80, 10).Value = 4001.3333
$ws.Cells.Item(80, 12).Value = 4001.3333
$ws.Cells.Item(80, 14).Value = -5997.3333
$ws.Cells.Item(83, 8).Value = 3658.75
$ws.Cells.Item(83, 10).Value = 4001.3333
$ws.Cells.Item(83, 12).Value = 20006.6665
$ws.Cells.Item(83, 14).Value = -29990.6665
$ws.Cells.Item(122, 8).Value = 2049.75
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 10).Value = 2099.5
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 12).Value = 6298.5
$ws.Cells.Item(122, 13).Value = -3550
$ws.Cells.Item(122, 14).Value = -11198.5
$ws.Cells.Item(126, 8).Value = 2573607.5
$ws.Cells.Item(126, 9).Value = 2927077.2
$ws.Cells.Item(126, 11).Value = 8781231.600000001
$ws.Cells.Item(126, 13).Value = -8778761.600000001
$ws.Cells.Item(132, 8).Value = 963569.4
$ws.Cells.Item(132, 9).Value = 1833420.4
$ws.Cells.Item(132, 11).Value = 5500261.199999999
$ws.Cells.Item(132, 13).Value = -5497731.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3321.3
$ws.Cells.Item(7, 9).Value = 3360.8
$ws.Cells.Item(7, 11).Value = 3360.8
$ws.Cells.Item(7, 13).Value = -3248.8
$ws.Cells.Item(22, 8).Value = 1511.4546
$ws.Cells.Item(22, 9).Value = 1186.75
$ws.Cells.Item(22, 11).Value = 1186.75
$ws.Cells.Item(22, 13).Value = -891.75
$ws.Cells.Item(27, 8).Value = 1511.4546
$ws.Cells.Item(27, 9).Value = 1186.75
$ws.Cells.Item(27, 11).Value = 1186.75
$ws.Cells.Item(27, 13).Value = -1079.75
$ws.Cells.Item(40, 8).Value = 10389.417
$ws.Cells.Item(40, 9).Value = 9453.223
$ws.Cells.Item(40, 11).Value = 9453.223
$ws.Cells.Item(40, 13).Value = -9317.223
$ws.Cells.Item(82, 8).Value = 1700.875
$ws.Cells.Item(82, 9).Value = 1700.875
$ws.Cells.Item(82, 11).Value = 1700.875
$ws.Cells.Item(82, 13).Value = -1339.875
$ws.Cells.Item(85, 8).Value = 1700.875
$ws.Cells.Item(85, 9).Value = 1700.875
$ws.Cells.Item(85, 11).Value = 1700.875
$ws.Cells.Item(85, 13).Value = -452.875
$ws.Cells.Item(126, 8).Value = 3321.3
$ws.Cells.Item(126, 9).Value = 3360.8
$ws.Cells.Item(126, 11).Value = 10082.4
$ws.Cells.Item(126, 13).Value = -7612.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 20577900
$ws.Cells.Item(136, 9).Value = 32681234
$ws.Cells.Item(136, 10).Value = 2233.4
$ws.Cells.Item(136, 11).Value = 98043702
$ws.Cells.Item(136, 12).Value = 6700.200000000001
$ws.Cells.Item(136, 13).Value = -98041152
$ws.Cells.Item(136, 14).Value = -11800.2
